$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set a cell value as literal text (avoids Excel auto-converting numeric-looking
# strings such as "99.35" or "42.894.91" into numbers), while keeping the cell
# unstyled (no explicit "Text" number format applied) to match the original file.
function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

Set-TextValue "D2" "42.894.91"
Set-TextValue "E2" "  -1.47%  "

Set-TextValue "D3" "2.302.03"
Set-TextValue "E3" "  -3.26%  "

Set-TextValue "E4" "  -0.04%  "

Set-TextValue "D5" "302.32"
Set-TextValue "E5" "  -2.46%  "

Set-TextValue "D6" "99.35"
Set-TextValue "E6" "  -5.50%  "

Set-TextValue "E7" "  -1.11%  "

Set-TextValue "E8" "  +0.06%  "

Set-TextValue "E9" "  -2.26%  "

Set-TextValue "D10" "34.85"
Set-TextValue "E10" "  -3.77%  "

Set-TextValue "B11" "OKB"
Set-TextValue "C11" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D11" "51.27"
Set-TextValue "E11" "  -4.05%  "

Set-TextValue "B12" "Dogecoin"
Set-TextValue "C12" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D12" "0.0796"
Set-TextValue "E12" "  -2.08%  "

Set-TextValue "E13" "  +0.34%  "

Set-TextValue "D14" "6.78"
Set-TextValue "E14" "  -3.39%  "

Set-TextValue "D15" "2.660.19"
Set-TextValue "E15" "  -3.20%  "

Set-TextValue "D16" "15.53"
Set-TextValue "E16" "  -0.91%  "

Set-TextValue "D17" "2.309.39"
Set-TextValue "E17" "  -2.82%  "

Set-TextValue "D18" "0.797"
Set-TextValue "E18" "  -1.98%  "

Set-TextValue "D19" "42.810.53"
Set-TextValue "E19" "  -1.55%  "

Set-TextValue "D20" "11.75"
Set-TextValue "E20" "  -2.28%  "

Set-TextValue "D21" "0.0₃0900"
Set-TextValue "E21" "  -2.33%  "

Set-TextValue "D22" "6.06"
Set-TextValue "E22" "  -3.81%  "

Set-TextValue "D23" "67.43"
Set-TextValue "E23" "  -1.51%  "

Set-TextValue "D24" "236.45"
Set-TextValue "E24" "  -2.26%  "

Set-TextValue "E25" "  -4.48%  "

Set-TextValue "E26" "  -3.92%  "

Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.16%  "

Set-TextValue "D28" "24.80"
Set-TextValue "E28" "  -3.90%  "

Set-TextValue "B29" "Toncoin"
Set-TextValue "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.24"
Set-TextValue "E29" "  -0.09%  "

Set-TextValue "B30" "InjectiveProtocol"
Set-TextValue "C30" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "34.44"
Set-TextValue "E30" "  -7.03%  "

Set-TextValue "D31" "164.69"
Set-TextValue "E31" "  +1.60%  "

Set-TextValue "E32" "  -4.19%  "

Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  -0.08%  "

Set-TextValue "D34" "5.03"
Set-TextValue "E34" "  -4.69%  "

Set-TextValue "E35" "  -4.89%  "

Set-TextValue "E36" "  -4.91%  "

Set-TextValue "E37" "  -5.22%  "

Set-TextValue "D38" "16.52"
Set-TextValue "E38" "  -10.38%  "

Set-TextValue "E39" "  -7.85%  "

Set-TextValue "E41" "  -4.26%  "

Set-TextValue "D42" "0.111"
Set-TextValue "E42" "  -2.82%  "

Set-TextValue "D43" "2.42"
Set-TextValue "E43" "  -9.22%  "

Set-TextValue "D44" "1.967.94"
Set-TextValue "E44" "  -3.03%  "

Set-TextValue "E45" "  -2.06%  "

Set-TextValue "D46" "18.35"
Set-TextValue "E46" "  -6.89%  "

Set-TextValue "D47" "9.82"
Set-TextValue "E47" "  -7.41%  "

Set-TextValue "E48" "  -8.56%  "

Set-TextValue "D49" "4.79"
Set-TextValue "E49" "  +0.77%  "

Set-TextValue "E50" "  -6.87%  "

Set-TextValue "D51" "2.531.06"
Set-TextValue "E51" "  -2.67%  "
